# Added a new skill, new alchemy item and new guide quests.
# Adds a new "Prisoners Escape" skill row (row 33) to the "Game Skills" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Prisoners Escape"
$ws.Range("C33").Value = 13
$ws.Range("D33").Value = "Escape the pits of torment and the delusions that haunt your mind from the time you spent in captivity. Escape with the fury and rage of a thousand men."
$ws.Range("E33").Value = 999
$ws.Range("F33").Value = 0.00175
$ws.Range("H33").Value = 0.0005
$ws.Range("I33").Value = 0.0005
$ws.Range("J33").Value = 0.0005
$ws.Range("N33").Value = 1
$ws.Range("P33").Value = 0.001
$ws.Range("Q33").Value = 9
$ws.Range("R33").Value = 0
